# Generate Report for Handback
# Updates the handoff/handback timestamps for the
# "e87b3ebd-18b1-49cf-b232-fe0371daea31.md" row (row 3) now that a new
# xliff round-trip has completed, and reflects the newest timestamp on
# the Overview sheet.

$wb = $excel.ActiveWorkbook

$newHandoffDatetimeZhCn  = "2016-08-12 18:59:45"
$newHandbackDatetimeZhCn = "2016-08-12 19:00:26"

$newHandoffDatetimeDeDe  = "2016-08-12 18:59:52"
$newHandbackDatetimeDeDe = "2016-08-12 19:00:36"

$newLatestHoXliffGenerateDate = "2016-08-12 18:59:52"

# zh-cn sheet: row 3 is the e87b3ebd-18b1-49cf-b232-fe0371daea31 file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = $newHandoffDatetimeZhCn
$wsZhCn.Range("K3").Value = $newHandbackDatetimeZhCn

# de-de sheet: row 3 is the e87b3ebd-18b1-49cf-b232-fe0371daea31 file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = $newHandoffDatetimeDeDe
$wsDeDe.Range("K3").Value = $newHandbackDatetimeDeDe

# Overview sheet: row 3 "Latest HO Xliff Generate Date" reflects the
# most recent handback across languages for that file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = $newLatestHoXliffGenerateDate
